$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), columns G..R -------------------------------
# Shared strings must land in the workbook's sst in the same order the
# authoring session produced them: G..P left-to-right, then R before Q
# (so FINAL.MATURITY gets the lower shared-string index than EXP.DATE).
$ws.Range("G1").Value = "INTEREST.RATE"
$ws.Range("H1").Value = "PROFIT.PAY.TERM"
$ws.Range("I1").Value = "INTEND.DATE"
$ws.Range("J1").Value = "CUST.REMARKS:1"
$ws.Range("K1").Value = "LIMIT.REFERENCE"
$ws.Range("L1").Value = "TAX.INTEREST.TYPE:1"
$ws.Range("M1").Value = "DRAWDOWN.ACCOUNT"
$ws.Range("N1").Value = "PRIN.LIQ.ACCT"
$ws.Range("O1").Value = "INT.LIQ.ACCT"
$ws.Range("P1").Value = "CHRG.LIQ.ACCT"
$ws.Range("R1").Value = "FINAL.MATURITY"
$ws.Range("Q1").Value = "EXP.DATE"

# --- Column widths ----------------------------------------------------------
# Existing columns C/D got narrower (no longer share D's generic "customWidth"
# only formatting), both now best-fit like the rest.
$ws.Columns.Item(3).ColumnWidth = 9.5
$ws.Columns.Item(4).ColumnWidth = 8.833333333333334

# Newly introduced columns, sized to fit their header text.
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 16
$ws.Columns.Item(9).ColumnWidth = 12
$ws.Columns.Item(10).ColumnWidth = 15.333333333333334
$ws.Columns.Item(11).ColumnWidth = 15.5
$ws.Columns.Item(12).ColumnWidth = 18.833333333333332
$ws.Columns.Item(13).ColumnWidth = 21.5
$ws.Columns.Item(14).ColumnWidth = 13.166666666666666
$ws.Columns.Item(15).ColumnWidth = 11.666666666666666
$ws.Columns.Item(16).ColumnWidth = 13.666666666666666

# --- View state: scrolled right with L13 the active selection --------------
$ws.Range("L13").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
